$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.136.28"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.656.01"
$ws.Range("E3").Value = "  -0.68%  "

$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("E6").Value = "  +0.94%  "

$ws.Range("E7").Value = "  -0.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2604"
$ws.Range("E8").Value = "  -2.31%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06339"
$ws.Range("E9").Value = "  +0.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.42"
$ws.Range("E10").Value = "  -2.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07775"
$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.678.06"
$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.498"
$ws.Range("E13").Value = "  +1.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5476"
$ws.Range("E14").Value = "  -0.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅8161"
$ws.Range("E15").Value = "  -1.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.49"
$ws.Range("E16").Value = "  +0.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.131.57"
$ws.Range("E17").Value = "  -0.69%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.570"
$ws.Range("E19").Value = "  -2.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.99"
$ws.Range("E20").Value = "  -0.89%  "

$ws.Range("E21").Value = "  -0.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.028"
$ws.Range("E22").Value = "  -0.70%  "

$ws.Range("E23").Value = "  -0.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "141.91"
$ws.Range("E24").Value = "  +1.67%  "

$ws.Range("E25").Value = "  +1.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.271"
$ws.Range("E26").Value = "  +0.95%  "

$ws.Range("E27").Value = "  +0.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.441"
$ws.Range("E28").Value = "  +1.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05941"
$ws.Range("E29").Value = "  -3.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.279"
$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.515"
$ws.Range("E31").Value = "  -2.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.248"
$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.576"
$ws.Range("E33").Value = "  -3.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9509"
$ws.Range("E34").Value = "  -2.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.792"
$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("E36").Value = "  -0.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5667"
$ws.Range("E37").Value = "  -1.64%  "

$ws.Range("E38").Value = "  +0.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.820"
$ws.Range("E39").Value = "  -3.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8476"
$ws.Range("E40").Value = "  -1.25%  "

$ws.Range("E41").Value = "  -0.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.67"
$ws.Range("E42").Value = "  +2.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.022.78"
$ws.Range("E43").Value = "  -0.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.800.29"
$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.19"
$ws.Range("E45").Value = "  -0.73%  "

$ws.Range("E46").Value = "  +0.28%  "

$ws.Range("E47").Value = "  +1.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.478"
$ws.Range("E48").Value = "  -0.64%  "

$ws.Range("E49").Value = "  -0.66%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.791"
$ws.Range("E50").Value = "  -3.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09708"
$ws.Range("E51").Value = "  -0.74%  "
